$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.436.77"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").Value = "2.274.03"
$ws.Range("E3").Value = "  -0.79%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "122.46"
$ws.Range("E5").Value = "  +5.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.25"
$ws.Range("E6").Value = "  -1.01%  "

$ws.Range("E7").Value = "  +2.34%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("E9").Value = "  +0.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.92"
$ws.Range("E10").Value = "  -2.59%  "

$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.25"
$ws.Range("E12").Value = "  +3.94%  "

$ws.Range("E13").Value = "  -1.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.42"
$ws.Range("E14").Value = "  -2.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.907"
$ws.Range("E15").Value = "  +2.95%  "

$ws.Range("D16").Value = "2.615.42"
$ws.Range("E16").Value = "  -0.86%  "

$ws.Range("D17").Value = "2.276.06"
$ws.Range("E17").Value = "  -0.22%  "

$ws.Range("D18").Value = "43.498.63"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("E19").Value = "  +0.84%  "

$ws.Range("E20").Value = "  -0.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.27"
$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.42"
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.31"
$ws.Range("E23").Value = "  +1.01%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.51"
$ws.Range("E24").Value = "  -3.92%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.89"
$ws.Range("E25").Value = "  -1.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.02"
$ws.Range("E26").Value = "  +2.96%  "

$ws.Range("E27").Value = "  +1.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.37"
$ws.Range("E28").Value = "  +0.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.36"
$ws.Range("E29").Value = "  -0.67%  "

$ws.Range("E30").Value = "  +0.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.67"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.69"
$ws.Range("E32").Value = "  +0.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0914"
$ws.Range("E33").Value = "  -1.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.72"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("E35").Value = "  +1.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.31"
$ws.Range("E36").Value = "  +13.29%  "

$ws.Range("E37").Value = "  +4.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.60"
$ws.Range("E38").Value = "  -2.02%  "

$ws.Range("E39").Value = "  -1.35%  "

$ws.Range("E40").Value = "  +5.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.98"
$ws.Range("E41").Value = "  -4.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.81"
$ws.Range("E42").Value = "  -0.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.238"
$ws.Range("E43").Value = "  -1.47%  "

$ws.Range("E44").Value = "  -0.14%  "

$ws.Range("E45").Value = "  -0.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.69"
$ws.Range("E46").Value = "  -11.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "73.86"
$ws.Range("E47").Value = "  +38.84%  "

$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.27"
$ws.Range("E48").Value = "  +0.12%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.57"
$ws.Range("E49").Value = "  -1.68%  "

$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.50"
$ws.Range("E51").Value = "  -1.34%  "
